# Applies numeric updates to H:N columns across multiple sheets
# as produced by the scheduled runner refresh of item price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 4847.5835
$ws.Range("I32").Value = 5825.5
$ws.Range("J32").Value = 3478.5
$ws.Range("K32").Value = 5825.5
$ws.Range("L32").Value = 3478.5
$ws.Range("M32").Value = -5499.5
$ws.Range("N32").Value = -4130.5

$ws.Range("H52").Value = 3077.8
$ws.Range("J52").Value = 2265.6667
$ws.Range("L52").Value = 6797.000100000001
$ws.Range("N52").Value = -7117.000100000001

$ws.Range("H55").Value = 220.57143
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 220.57143
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 220.57143
$ws.Range("M55").Value = ""
$ws.Range("N55").Value = -648.57143

$ws.Range("H62").Value = 2002227.5
$ws.Range("I62").Value = 2002227.5
$ws.Range("K62").Value = 2002227.5
$ws.Range("M62").Value = -2001603.5

$ws.Range("H65").Value = 2002227.5
$ws.Range("I65").Value = 2002227.5
$ws.Range("K65").Value = 10011137.5
$ws.Range("M65").Value = -10008017.5

$ws.Range("H76").Value = 7525.6787
$ws.Range("I76").Value = 7389.8667
$ws.Range("K76").Value = 7389.8667
$ws.Range("M76").Value = -7074.8667

$ws.Range("H79").Value = 7525.6787
$ws.Range("I79").Value = 7389.8667
$ws.Range("K79").Value = 7389.8667
$ws.Range("M79").Value = -6297.8667

$ws.Range("H98").Value = 1087.1875
$ws.Range("I98").Value = 1093.0667
$ws.Range("K98").Value = 1093.0667
$ws.Range("M98").Value = 404.9332999999999

$ws.Range("H100").Value = 4749.5
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = ""

$ws.Range("H107").Value = 1612.0625
$ws.Range("I107").Value = 1519.6154
$ws.Range("K107").Value = 1519.6154
$ws.Range("M107").Value = 400.3846000000001

$ws.Range("H122").Value = 1087.1875
$ws.Range("I122").Value = 1093.0667
$ws.Range("K122").Value = 3279.2001
$ws.Range("M122").Value = -829.2001

$ws.Range("H132").Value = 3276.2812
$ws.Range("I132").Value = 3257.3
$ws.Range("K132").Value = 9771.900000000001
$ws.Range("M132").Value = -7241.900000000001

$ws.Range("H135").Value = 7466.8237
$ws.Range("I135").Value = 1516.9
$ws.Range("K135").Value = 13652.1
$ws.Range("M135").Value = -11117.1

$ws.Range("H138").Value = 3328.5957
$ws.Range("I138").Value = 1580.25
$ws.Range("J138").Value = 4623.6665
$ws.Range("K138").Value = 4740.75
$ws.Range("L138").Value = 13870.9995
$ws.Range("M138").Value = 399.25
$ws.Range("N138").Value = -24150.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2511
$ws.Range("I5").Value = 2076.625
$ws.Range("K5").Value = 2076.625
$ws.Range("M5").Value = -1964.625

$ws.Range("H61").Value = 5183.294
$ws.Range("I61").Value = 4994.393
$ws.Range("K61").Value = 4994.393
$ws.Range("M61").Value = -4782.393

$ws.Range("H122").Value = 1682.0358
$ws.Range("I122").Value = 1539.5217
$ws.Range("K122").Value = 4618.5651
$ws.Range("M122").Value = -2168.5651

$ws.Range("H132").Value = 4533.607
$ws.Range("I132").Value = 4601.75
$ws.Range("K132").Value = 13805.25
$ws.Range("M132").Value = -11275.25

$ws.Range("H136").Value = 5183.294
$ws.Range("I136").Value = 4994.393
$ws.Range("K136").Value = 14983.179
$ws.Range("M136").Value = -12433.179

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2511
$ws.Range("I4").Value = 2076.625
$ws.Range("K4").Value = 2076.625
$ws.Range("M4").Value = -1961.625

$ws.Range("H80").Value = 1292.5
$ws.Range("J80").Value = 966.44446
$ws.Range("L80").Value = 966.44446
$ws.Range("N80").Value = -2962.44446

$ws.Range("H83").Value = 1292.5
$ws.Range("J83").Value = 966.44446
$ws.Range("L83").Value = 4832.2223
$ws.Range("N83").Value = -14816.2223

$ws.Range("H105").Value = 1801.7273
$ws.Range("I105").Value = 1757.8889
$ws.Range("K105").Value = 1757.8889
$ws.Range("M105").Value = -10.88889999999992

$ws.Range("H134").Value = 2881.375
$ws.Range("I134").Value = 2864.6428
$ws.Range("K134").Value = 8593.928400000001
$ws.Range("M134").Value = -6058.928400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7487.522
$ws.Range("I31").Value = 8195.799999999999
$ws.Range("K31").Value = 8195.799999999999
$ws.Range("M31").Value = -7900.799999999999

$ws.Range("H34").Value = 7487.522
$ws.Range("I34").Value = 8195.799999999999
$ws.Range("K34").Value = 8195.799999999999
$ws.Range("M34").Value = -7993.799999999999

$ws.Range("H87").Value = 30000
$ws.Range("I87").Value = 30000
$ws.Range("K87").Value = 30000
$ws.Range("M87").Value = -28814

$ws.Range("H90").Value = 30000
$ws.Range("I90").Value = 30000
$ws.Range("K90").Value = 90000
$ws.Range("M90").Value = -84072

$ws.Range("H107").Value = 1185.7916
$ws.Range("I107").Value = 1184.5454
$ws.Range("J107").Value = 1199.5
$ws.Range("K107").Value = 1184.5454
$ws.Range("L107").Value = 1199.5
$ws.Range("M107").Value = 735.4546
$ws.Range("N107").Value = -5039.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4347.75
$ws.Range("I70").Value = 6000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5730

$ws.Range("H73").Value = 4347.75
$ws.Range("I73").Value = 6000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -5064

$ws.Range("H102").Value = 1570.2
$ws.Range("I102").Value = 1489.1111
$ws.Range("K102").Value = 1489.1111
$ws.Range("M102").Value = 132.8888999999999

$ws.Range("H108").Value = 70310
$ws.Range("I108").Value = 40621
$ws.Range("K108").Value = 40621
$ws.Range("M108").Value = -36781

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2298.5173
$ws.Range("I22").Value = 2132.8572
$ws.Range("J22").Value = 2453.1333
$ws.Range("K22").Value = 2132.8572
$ws.Range("L22").Value = 2453.1333
$ws.Range("M22").Value = -1837.8572
$ws.Range("N22").Value = -3043.1333

$ws.Range("H27").Value = 2298.5173
$ws.Range("I27").Value = 2132.8572
$ws.Range("J27").Value = 2453.1333
$ws.Range("K27").Value = 2132.8572
$ws.Range("L27").Value = 2453.1333
$ws.Range("M27").Value = -2025.8572
$ws.Range("N27").Value = -2667.1333

$ws.Range("H40").Value = 10049.115
$ws.Range("I40").Value = 6989.7896
$ws.Range("K40").Value = 6989.7896
$ws.Range("M40").Value = -6853.7896

$ws.Range("H93").Value = 1151.9231
$ws.Range("I93").Value = 1151.9231
$ws.Range("K93").Value = 1151.9231
$ws.Range("M93").Value = 96.07690000000002

$ws.Range("H133").Value = 89999
$ws.Range("J133").Value = 89999
$ws.Range("L133").Value = 89999
$ws.Range("N133").Value = -95059

$ws.Range("H136").Value = 2250
$ws.Range("I136").Value = 2250
$ws.Range("K136").Value = 6750
$ws.Range("M136").Value = -4200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4075.8333
$ws.Range("I132").Value = 4100.125
$ws.Range("K132").Value = 12300.375
$ws.Range("M132").Value = -9770.375
